$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ip_address_list")

$ws.Range("D1").Value = "adfadf"
$ws.Range("D2").Value = "FortiClient Edcha Ex2dp78kxp30`nsfg`nsffgs`ngsfg`ndfa`nf`nsfg`nsfg"
$ws.Range("D3").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.245`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK`nsf"
$ws.Range("D4").Value = "Teleflex `nsfg`nsgf"
$ws.Range("D5").Value = "PC:192.168.14.240`nf`nf"
$ws.Range("D7").Value = "Kamera VS-S160MX :192.168.0.186`nddfaajhdf"
$ws.Range("D9").Value = "afd"
$ws.Range("D10").Value = "PC:`t10.96.205.175`najpodkjfa"
$ws.Range("D11").Value = "XG-X2900:`t`t10.101.28.175`nadfddd`nadfdddd`nadfdd`nadfg`nfg`nf`ndf`nadf`nf`nadf"

$settings = $wb.Worksheets.Item("Settings")
$settings.Range("B3").Value = 1
